# Automatic update of files.
# Applies the row-level corrections described by the upstream diff:
#  - Several rows only had their "B" (Taxonsorteringsordning) value bumped by 1.
#  - Rows 16/17/18, 21/22/23 and 30/31 had their full record content rotated
#    between rows (A, B, D, E, F, G, H, Q, R and, where relevant, K/L/M/N/AC).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Group 1: simple increments of column B only
# ---------------------------------------------------------------------------
$ws.Range("B6").Value  = 83217
$ws.Range("B7").Value  = 80385
$ws.Range("B9").Value  = 83091
$ws.Range("B11").Value = 91773
$ws.Range("B12").Value = 83208
$ws.Range("B13").Value = 83208
$ws.Range("B14").Value = 83217
$ws.Range("B24").Value = 83091
$ws.Range("B25").Value = 91773
$ws.Range("B26").Value = 91773
$ws.Range("B27").Value = 83217
$ws.Range("B28").Value = 83217
$ws.Range("B29").Value = 92229
$ws.Range("B32").Value = 91810
$ws.Range("B33").Value = 83217
$ws.Range("B34").Value = 83091

# ---------------------------------------------------------------------------
# Group 2: rows 16, 17, 18 - record content rotates (row16<-row18, row17<-row16,
# row18<-row17), each with its B value incremented by 1.
# ---------------------------------------------------------------------------
$ws.Range("A16").Value = 131066761
$ws.Range("B16").Value = 91773
$ws.Range("D16").Value = "LC"
$ws.Range("E16").Value = 5447
$ws.Range("F16").Value = "Vedticka"
$ws.Range("G16").Value = "Fuscoporia viticola"
$ws.Range("H16").Value = "(Schwein.) Murrill"
$ws.Range("Q16").Value = 425072
$ws.Range("R16").Value = 6712273

$ws.Range("A17").Value = 131066782
$ws.Range("B17").Value = 91824
$ws.Range("D17").Value = "NT"
$ws.Range("E17").Value = 1204
$ws.Range("F17").Value = "Gränsticka"
$ws.Range("G17").Value = "Phellopilus nigrolimitatus"
$ws.Range("H17").Value = "(Romell) Niemelä, T.Wagner & M.Fisch."
$ws.Range("Q17").Value = 425059
$ws.Range("R17").Value = 6712253

$ws.Range("A18").Value = 131066768
$ws.Range("B18").Value = 91810
$ws.Range("D18").Value = "NT"
$ws.Range("E18").Value = 1202
$ws.Range("F18").Value = "Ullticka"
$ws.Range("G18").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H18").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("Q18").Value = 425256
$ws.Range("R18").Value = 6712203

# ---------------------------------------------------------------------------
# Group 3: rows 21, 22, 23 - record content rotates (row21<-row22, row22<-row23,
# row23<-row21). Row 21 also loses its K/L/M/N/AC cells while row 23 gains them
# (the "äldre spår" / "Ringhack på gran" woodpecker-track annotation moves from
# row 21 to row 23 along with the rest of that record).
# ---------------------------------------------------------------------------
$ws.Range("A21").Value = 131066766
$ws.Range("B21").Value = 92181
$ws.Range("D21").Value = "VU"
$ws.Range("E21").Value = 2062
$ws.Range("F21").Value = "Ulltickeporing"
$ws.Range("G21").Value = "Skeletocutis brevispora"
$ws.Range("H21").Value = "Niemelä"
$ws.Range("K21").ClearContents()
$ws.Range("L21").ClearContents()
$ws.Range("M21").ClearContents()
$ws.Range("N21").ClearContents()
$ws.Range("Q21").Value = 425069
$ws.Range("R21").Value = 6712285
$ws.Range("AC21").ClearContents()

$ws.Range("A22").Value = 131066778
$ws.Range("B22").Value = 81230
$ws.Range("D22").Value = "NT"
$ws.Range("E22").Value = 1049
$ws.Range("F22").Value = "Kortskaftad ärgspik"
$ws.Range("G22").Value = "Microcalicium ahlneri"
$ws.Range("H22").Value = "Tibell"
$ws.Range("Q22").Value = 425336
$ws.Range("R22").Value = 6712202

$ws.Range("A23").Value = 131066774
$ws.Range("B23").Value = 57884
$ws.Range("D23").Value = "NT"
$ws.Range("E23").Value = 100109
$ws.Range("F23").Value = "Tretåig hackspett"
$ws.Range("G23").Value = "Picoides tridactylus"
$ws.Range("H23").Value = "(Linnaeus, 1758)"
$ws.Range("M23").Value = "äldre spår"
$ws.Range("Q23").Value = 425250
$ws.Range("R23").Value = 6712265
$ws.Range("AC23").Value = "Ringhack på gran"

# ---------------------------------------------------------------------------
# Group 4: rows 30, 31 swap record content (row30<-row31 with B+1, row31<-row30
# unchanged B). The woodpecker-track annotation (K/L/M/N/AC) moves from row 30
# to row 31 along with the rest of that record.
# ---------------------------------------------------------------------------
$ws.Range("A30").Value = 131066776
$ws.Range("B30").Value = 80351
$ws.Range("E30").Value = 2081
$ws.Range("F30").Value = "Skrovellav"
$ws.Range("G30").Value = "Lobaria scrobiculata"
$ws.Range("H30").Value = "(Scop.) DC."
$ws.Range("K30").ClearContents()
$ws.Range("L30").ClearContents()
$ws.Range("M30").ClearContents()
$ws.Range("N30").ClearContents()
$ws.Range("Q30").Value = 425069
$ws.Range("R30").Value = 6712285
$ws.Range("AC30").ClearContents()

$ws.Range("A31").Value = 131066772
$ws.Range("B31").Value = 57884
$ws.Range("E31").Value = 100109
$ws.Range("F31").Value = "Tretåig hackspett"
$ws.Range("G31").Value = "Picoides tridactylus"
$ws.Range("H31").Value = "(Linnaeus, 1758)"
$ws.Range("M31").Value = "äldre spår"
$ws.Range("Q31").Value = 425301
$ws.Range("R31").Value = 6712219
$ws.Range("AC31").Value = "Ringhack på gran"
